{"js": "// The site footer used to show, after the last \"Requisitos\" entry\n// (\"LOB1012: Estat\u00edstica (Requisito)\"), an empty paragraph followed by\n// two more paragraphs:\n//   \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n//    pages. Original theme under Creative Commons Attribution\"\n// This rebuild of the site dropped that footer block, so remove those\n// three paragraphs while leaving the \"LOB1012\u2026\" paragraph (and the\n// blank paragraph that originally followed the footer) untouched.\n\nconst body = context.document.body;\nconst results = body.search(\"LOB1012\", { matchCase: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const anchorPara = results.items[0].paragraphs.getFirst();\n  const blank = anchorPara.getNext();\n  const jupiterPara = blank.getNext();\n  const copyrightPara = jupiterPara.getNext();\n\n  blank.load(\"text\");\n  jupiterPara.load(\"text\");\n  copyrightPara.load(\"text\");\n  await context.sync();\n\n  // Only delete if the paragraphs look like what we expect, so the\n  // script is a no-op (rather than destructive) if the document has\n  // already changed shape.\n  if (blank.text === \"\" &&\n      jupiterPara.text === \"Ver no Jupiter Salvar em pdf Salvar em docx\" &&\n      copyrightPara.text.indexOf(\"\u00a9 2020\") === 0) {\n    blank.delete();\n    jupiterPara.delete();\n    copyrightPara.delete();\n    await context.sync();\n  }\n}\n", "ps1": "# The site footer used to show, after the last \"Requisitos\" entry\n# (\"LOB1012: Estat\u00edstica (Requisito)\"), an empty paragraph followed by\n# two more paragraphs:\n#   \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n#    pages. Original theme under Creative Commons Attribution\"\n# This rebuild of the site dropped that footer block, so remove those\n# three paragraphs while leaving the \"LOB1012\u2026\" paragraph (and the\n# blank paragraph that originally followed the footer) untouched.\n\n$d = $word.ActiveDocument\n\n$search = $d.Content\n$found = $search.Find.Execute(\"LOB1012\")\n\nif ($found) {\n    $anchorPara = $search.Paragraphs(1)\n    $blank = $anchorPara.Next()\n    $jupiterPara = $blank.Next()\n    $copyrightPara = $jupiterPara.Next()\n\n    if ($blank.Range.Text -eq \"`r\" -and\n        $jupiterPara.Range.Text -eq \"Ver no Jupiter Salvar em pdf Salvar em docx`r\" -and\n        $copyrightPara.Range.Text.Contains(\"2020 . Contact: luizeleno@usp.br\")) {\n        $delRange = $d.Range($blank.Range.Start, $copyrightPara.Range.End)\n        $delRange.Delete()\n    }\n}\n"}
